# Updated cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
  # Force the cell to be treated as text so numeric-looking strings
  # (e.g. "1.000", "13.60", "0.000007577") are preserved exactly,
  # then restore the default "Normal" style so no stray number format
  # is left behind on the cell.
  $cell.NumberFormat = "@"
  $cell.Value = $value
  $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '30.451.13'
$ws.Cells.Item(2, 5).Value = '  +1.02%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.877.73'
Set-TextValue $ws.Cells.Item(4, 4) '1.000'
Set-TextValue $ws.Cells.Item(5, 4) '247.11'
$ws.Cells.Item(5, 5).Value = '  +5.69%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.15%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.4772'
$ws.Cells.Item(7, 5).Value = '  +1.98%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.2897'
$ws.Cells.Item(8, 5).Value = '  +1.84%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.06520'
$ws.Cells.Item(9, 5).Value = '  +0.91%  '
Set-TextValue $ws.Cells.Item(10, 4) '21.89'
$ws.Cells.Item(10, 5).Value = '  +4.52%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.07725'
$ws.Cells.Item(11, 5).Value = '  -0.33%  '
Set-TextValue $ws.Cells.Item(12, 4) '97.11'
$ws.Cells.Item(12, 5).Value = '  +4.03%  '
Set-TextValue $ws.Cells.Item(13, 4) '0.7391'
$ws.Cells.Item(13, 5).Value = '  +8.95%  '
Set-TextValue $ws.Cells.Item(14, 4) '1.877.13'
$ws.Cells.Item(14, 5).Value = '  +0.56%  '
Set-TextValue $ws.Cells.Item(16, 4) '273.39'
$ws.Cells.Item(16, 5).Value = '  +2.89%  '
Set-TextValue $ws.Cells.Item(17, 4) '30.445.65'
$ws.Cells.Item(17, 5).Value = '  +1.06%  '
Set-TextValue $ws.Cells.Item(18, 4) '13.60'
Set-TextValue $ws.Cells.Item(19, 4) '0.000007577'
$ws.Cells.Item(19, 5).Value = '  +0.09%  '
Set-TextValue $ws.Cells.Item(20, 4) '1.000'
$ws.Cells.Item(20, 5).Value = '  -0.11%  '
Set-TextValue $ws.Cells.Item(21, 4) '2.123.71'
$ws.Cells.Item(21, 5).Value = '  +0.30%  '
$ws.Cells.Item(22, 5).Value = '  -0.13%  '
Set-TextValue $ws.Cells.Item(23, 4) '5.257'
$ws.Cells.Item(23, 5).Value = '  +2.67%  '
Set-TextValue $ws.Cells.Item(24, 4) '6.182'
$ws.Cells.Item(24, 5).Value = '  +1.57%  '
Set-TextValue $ws.Cells.Item(25, 4) '9.332'
$ws.Cells.Item(25, 5).Value = '  +0.25%  '
Set-TextValue $ws.Cells.Item(26, 4) '163.77'
$ws.Cells.Item(26, 5).Value = '  -0.79%  '
Set-TextValue $ws.Cells.Item(27, 4) '18.85'
$ws.Cells.Item(27, 5).Value = '  +2.18%  '
Set-TextValue $ws.Cells.Item(28, 4) '1.943'
$ws.Cells.Item(28, 5).Value = '  +3.52%  '
$ws.Cells.Item(29, 5).Value = '  +0.74%  '
Set-TextValue $ws.Cells.Item(30, 4) '0.09957'
$ws.Cells.Item(30, 5).Value = '  +0.53%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.520'
$ws.Cells.Item(31, 5).Value = '  +4.86%  '
$ws.Cells.Item(32, 5).Value = '  +2.65%  '
Set-TextValue $ws.Cells.Item(33, 4) '4.066'
$ws.Cells.Item(33, 5).Value = '  +2.35%  '
Set-TextValue $ws.Cells.Item(34, 4) '0.04791'
$ws.Cells.Item(34, 5).Value = '  +2.99%  '
Set-TextValue $ws.Cells.Item(35, 4) '1.125'
$ws.Cells.Item(35, 5).Value = '  +1.19%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.7007'
$ws.Cells.Item(36, 5).Value = '  +2.23%  '
Set-TextValue $ws.Cells.Item(37, 4) '2.715'
$ws.Cells.Item(37, 5).Value = '  +0.03%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.01872'
$ws.Cells.Item(38, 5).Value = '  +2.41%  '
Set-TextValue $ws.Cells.Item(39, 4) '2.728'
$ws.Cells.Item(39, 5).Value = '  -0.84%  '
Set-TextValue $ws.Cells.Item(40, 4) '6.328'
$ws.Cells.Item(40, 5).Value = '  +0.87%  '
Set-TextValue $ws.Cells.Item(41, 4) '71.10'
$ws.Cells.Item(41, 5).Value = '  +0.23%  '
Set-TextValue $ws.Cells.Item(42, 4) '1.951'
$ws.Cells.Item(42, 5).Value = '  +3.92%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.4209'
$ws.Cells.Item(43, 5).Value = '  +4.43%  '
$ws.Cells.Item(44, 5).Value = '  -0.10%  '
Set-TextValue $ws.Cells.Item(45, 4) '0.8370'
$ws.Cells.Item(45, 5).Value = '  +0.81%  '
Set-TextValue $ws.Cells.Item(46, 4) '102.89'
$ws.Cells.Item(46, 5).Value = '  +0.92%  '
Set-TextValue $ws.Cells.Item(47, 4) '9.261'
$ws.Cells.Item(47, 5).Value = '  +1.93%  '
Set-TextValue $ws.Cells.Item(48, 4) '7.081'
Set-TextValue $ws.Cells.Item(49, 4) '35.63'
$ws.Cells.Item(49, 5).Value = '  +5.07%  '
Set-TextValue $ws.Cells.Item(50, 4) '925.02'
$ws.Cells.Item(50, 5).Value = '  -0.16%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.05648'
$ws.Cells.Item(51, 5).Value = '  +1.36%  '
